$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1 (Amalzar_Madhyamik_24-7): move the selection from B5 to F2
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Checklist" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)   # Madhyamik_Follow up_26-9 - source for text style
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$checklist = $wb.Worksheets.Add($null, $lastSheet)
$checklist.Name = "Checklist"

# Fill column B (item names) first, starting with row 2, then the header
# row cells - this reproduces the exact shared-string insertion order of
# the source workbook.
$checklist.Range("B2").Value = "OPV 0"
$checklist.Range("C1").Value = "Completion Date"
$checklist.Range("D1").Value = "Base Date"
$checklist.Range("B1").Value = "Item Name"
$checklist.Range("B3").Value = "Pentavalent 3"
$checklist.Range("B4").Value = "Pentavalent 2"
$checklist.Range("B5").Value = "Pentavalent 1"
$checklist.Range("B6").Value = "Albendazole 5"
$checklist.Range("B7").Value = "Albendazole 4"
$checklist.Range("B8").Value = "Albendazole 3"
$checklist.Range("B9").Value = "Albendazole 2"
$checklist.Range("B10").Value = "Albendazole 1"
$checklist.Range("B11").Value = "Vitamin A 5"
$checklist.Range("B12").Value = "Vitamin A 3"
$checklist.Range("B13").Value = "Vitamin A 2"
$checklist.Range("B14").Value = "Vitamin A 1"
$checklist.Range("B15").Value = "Measles 2"

$checklist.Range("A1").Value = "Enrolment UUID"
for ($r = 2; $r -le 15; $r++) {
    $checklist.Range("A" + $r).Value = "81286f81-e70c-4428-ad63-589a2d36e0f6"
}

$checkDate = (Get-Date -Year 2003 -Month 6 -Day 6).Date
for ($r = 2; $r -le 15; $r++) {
    $checklist.Range("C" + $r).Value = $checkDate
    $checklist.Range("D" + $r).Value = $checkDate
}

# Match formatting/styles of the sibling "Madhyamik_Follow up_26-9" sheet:
# text cells -> style of B2 (plain Calibri), date cells -> style of F2 on
# the first sheet (numFmtId 14 date format, right aligned).
$textStyleSrc = $ws3.Range("B2")
$textStyleSrc.Copy()
$checklist.Range("A1:D1").PasteSpecial(-4122)
$checklist.Range("A2:B15").PasteSpecial(-4122)

$dateStyleSrc = $ws1.Range("F2")
$dateStyleSrc.Copy()
$checklist.Range("C2:D15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$checklist.Range("B16").Select()

# ---------------------------------------------------------------------------
# 3. Row 1 height on "ambos 3-10" grows from 78 to 83
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("ambos 3-10")
$ws5.Rows.Item(1).RowHeight = 83
